# compelete missile image effect #4
#
# Updates FrameCount (F) / FrameTime (G) values for several missile rows
# on the "Missile" sheet, and moves the sheet's scroll position / active
# selection to match the author's final view (topLeftCell A10, selection
# G24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 (Id 13): FrameCount 1 -> 4, FrameTime 1 -> 2
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 2

# Row 20 (Id 16): FrameCount 1 -> 2, FrameTime 1 -> 2
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 2

# Row 21 (Id 17): FrameCount 1 -> 4, FrameTime 1 -> 2
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 2

# Row 22 (Id 18): FrameCount 1 -> 2, FrameTime 1 -> 3
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 3

# Row 25 (Id 21): FrameCount 1 -> 4, FrameTime 1 -> 3
$ws.Range("F25").Value = 4
$ws.Range("G25").Value = 3

# Scroll the view down (topLeftCell A13 -> A10) and move the active
# selection from D29 to G24.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G24").Select() | Out-Null
